$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update H38 (gun power for ammo_12x76_zhekan / Perf slugs): 2.7 -> 2.5
$ws.Range("H38").Value = 2.5

# Update H39 formula (gun power for ammo_12x70_buck / DMG buckshot): 9*0.42 -> 9*0.4
$ws.Range("H39").Formula = "=9*0.4"

# Update the active selection to match the reverted commit's cursor position
$ws.Range("N22").Select()
